{"js": "// Apply the five Zulu wording corrections described in the diff.\n// Each edit is a small, localized text substitution inside a single\n// run, so we locate the run via body.search() (unique surrounding\n// text) and replace just the changed fragment with insertText using\n// the \"Replace\" insert location, which preserves the run/paragraph\n// formatting of the matched range.\n\nconst edits = [\n  {\n    find: \"i-email ithimbeni locwaningo ku-\",\n    replace: \"i-email ithimba locwaningo ku-\",\n  },\n  {\n    find: \"ulwazi lwakho luyimfihlo futhi kuqinisekise ukuthi uzizwa ukhululekile\",\n    replace: \"ulwazi lwakho luyimfihlo futhi kuqinisekise ukuthi uzizwe ukhululekile\",\n  },\n  {\n    find: \"Hare), kanti iMenenja yocwaningo nguZamakhanya Makhanya\",\n    replace: \"Hare), kanye neMenenja yocwaningo uZamakhanya Makhanya\",\n  },\n  {\n    find: \"Uma unemibuzo noma okukukhathazayo mayelana namalungelo\",\n    replace: \"Uma unemibuzo noma kukhona okukukhathazayo mayelana namalungelo\",\n  },\n  {\n    find: \"Uma ufunde futhi waqonda idokhumenti engenhla\",\n    replace: \"Uma ufunde futhi waqonda incwadi engenhla\",\n  },\n];\n\nfor (const edit of edits) {\n  const results = context.document.body.search(edit.find, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + edit.find);\n  }\n\n  for (const result of results.items) {\n    result.insertText(edit.replace, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply the five Zulu wording corrections described in the diff.\n# Each edit is a small, localized text substitution inside a single\n# run. We use Find/Execute with a short, unique surrounding snippet\n# (so only the intended occurrence is matched) and Replace:=2\n# (wdReplaceAll, though only one match exists) to swap in the\n# corrected wording. MatchWholeWord is left False because several\n# search strings end in punctuation (e.g. a trailing hyphen) that\n# would never satisfy a whole-word boundary.\n\n$d = $word.ActiveDocument\n\n$edits = @(\n    @{ Find = \"i-email ithimbeni locwaningo ku-\"; Replace = \"i-email ithimba locwaningo ku-\" },\n    @{ Find = \"ulwazi lwakho luyimfihlo futhi kuqinisekise ukuthi uzizwa ukhululekile\"; Replace = \"ulwazi lwakho luyimfihlo futhi kuqinisekise ukuthi uzizwe ukhululekile\" },\n    @{ Find = \"Hare), kanti iMenenja yocwaningo nguZamakhanya Makhanya\"; Replace = \"Hare), kanye neMenenja yocwaningo uZamakhanya Makhanya\" },\n    @{ Find = \"Uma unemibuzo noma okukukhathazayo mayelana namalungelo\"; Replace = \"Uma unemibuzo noma kukhona okukukhathazayo mayelana namalungelo\" },\n    @{ Find = \"Uma ufunde futhi waqonda idokhumenti engenhla\"; Replace = \"Uma ufunde futhi waqonda incwadi engenhla\" }\n)\n\nforeach ($edit in $edits) {\n    $rng = $d.Content\n    $found = $rng.Find.Execute($edit.Find, $false, $false, $false, $false, $false, $true, 0, $false, $edit.Replace, 2)\n    if (-not $found) {\n        throw \"Search text not found: $($edit.Find)\"\n    }\n}\n"}
